$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that sits between the "D" run and
#    the "elete lehetoseg(most " run in the "Megoldani azt..." bullet.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2) The bullet "IDataError interface implementálása a UserProfile(Task) es
#    MyProfile-hoz" is replaced by the text/runs of the following bullet
#    ("IDataError kezdetit piros keret megoldasa (nem sikerult meg
#    megoldanom)"), and that following bullet is then removed - the two
#    bullets collapse into a single paragraph.
# ---------------------------------------------------------------------------
$oldPara = 0
$newPara = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text
    if ($txt.StartsWith("IDataError interface")) {
        $oldPara = $i
    } elseif ($txt.StartsWith("IDataError kezdetit")) {
        $newPara = $i
    }
}

$oldRange = $d.Paragraphs($oldPara).Range
$newRange = $d.Paragraphs($newPara).Range

# Ranges covering just the paragraph's content, without its trailing
# paragraph mark, so the formatted runs can be copied across untouched.
$oldContent = $d.Range($oldRange.Start, $oldRange.End - 1)
$newContent = $d.Range($newRange.Start, $newRange.End - 1)

$newContent.Copy()
$oldContent.Paste()

# The source bullet is now a duplicate - delete it (including its paragraph
# mark) so the two bullets merge into a single paragraph.
$d.Paragraphs($newPara).Range.Delete()

# ---------------------------------------------------------------------------
# 3) Inside the "Create utan listat frissit es akar odascrollozik ..."
#    bullet, split the run text "odascrollozik" into "odasc" / "rollozik"
#    and re-insert the (now orphaned) "_GoBack" bookmark at that split
#    point - exactly mirroring where it used to live before step 1.
# ---------------------------------------------------------------------------
$hit = $d.Content
$found = $hit.Find.Execute("odascrollozik", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

# "odascrollozik" -> "odasc" | "rollozik" : the split sits 5 characters in.
$splitPos = $hit.Start + 5
$splitRange = $d.Range($splitPos, $splitPos)
$d.Bookmarks.Add("_GoBack", $splitRange)
